# Added 10x versions 2 & 3 to assay_type in HCA scrnaseq
#
# Splits the old single "scRNAseq-10xGenomics" entry in the
# 'assay_type list' sheet into two explicit, versioned entries:
#   scRNAseq-10xGenomics-v2
#   scRNAseq-10xGenomics-v3
# and updates the assay_type data validation on the main sheet so it
# points at the now-6-row list instead of the old 5-row list.

$wb = $excel.ActiveWorkbook

# --- 'assay_type list' sheet: shift existing rows down one slot, then
#     write the two new 10x-version rows at the top. -----------------
$ws3 = $wb.Worksheets.Item("assay_type list")

$ws3.Range("A6").Value = $ws3.Range("A5").Value()   # SNARE2-RNAseq
$ws3.Range("A5").Value = $ws3.Range("A4").Value()   # snRNAseq
$ws3.Range("A4").Value = $ws3.Range("A3").Value()   # sciRNAseq
$ws3.Range("A3").Value = $ws3.Range("A2").Value()   # scRNAseq
$ws3.Range("A2").Value = "scRNAseq-10xGenomics-v3"
$ws3.Range("A1").Value = "scRNAseq-10xGenomics-v2"

# --- 'Export as TSV' sheet: point the assay_type (column J) list
#     validation at the new $A$1:$A$6 range and generalize the error
#     message so it no longer enumerates the (now longer) list. ------
$ws1 = $wb.Worksheets.Item("Export as TSV")
$colJ = $ws1.Range("J2:J1048576")
$colJ.Validation.Modify(3, 1, 3, "='assay_type list'!`$A`$1:`$A`$6")
$colJ.Validation.ErrorMessage = "Value must come from assay_type list."
